$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source row to replicate formatting/values from (last existing data row)
$srcRow = 225
$srcRange = $ws.Range("A" + $srcRow + ":J" + $srcRow)

# Date serial for A225 is 45781; new rows continue the daily sequence.
$startSerial = 45782
$firstNewRow = 226
$lastNewRow = 231

for ($row = $firstNewRow; $row -le $lastNewRow; $row++) {
    $destRange = $ws.Range("A" + $row + ":J" + $row)
    $srcRange.Copy($destRange) | Out-Null
    $ws.Cells.Item($row, 1).Value = $startSerial + ($row - $firstNewRow)
}

Write-Output "Added rows $firstNewRow to $lastNewRow"
